$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 35
$ws1.Range("F5").Value = 5425
$ws1.Range("F6").Value = 5425
$ws1.Range("F7").Value = 191
$ws1.Range("F9").Value = 540
$ws1.Range("F11").Value = 1216
$ws1.Range("F14").Value = 779
$ws1.Range("F15").Value = 6387
$ws1.Range("F16").Value = 36
$ws1.Range("F18").Value = 115
$ws1.Range("F19").Value = 3812
$ws1.Range("F20").Value = 274
$ws1.Range("F23").Value = 4097
$ws1.Range("F24").Value = 4010
$ws1.Range("F26").Value = 194
$ws1.Range("F27").Value = 256
$ws1.Range("F34").Value = 35
$ws1.Range("F35").Value = 7170
$ws1.Range("F36").Value = 36
$ws1.Range("F37").Value = 1186
$ws1.Range("F38").Value = 582
$ws1.Range("F39").Value = 111
$ws1.Range("F40").Value = 976
$ws1.Range("F42").Value = 1446
$ws1.Range("F44").Value = 781
$ws1.Range("F46").Value = 3331
$ws1.Range("F49").Value = 801

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 4
$ws2.Range("F11").Value = 68
$ws2.Range("F17").Value = 13
$ws2.Range("F18").Value = 147
$ws2.Range("F23").Value = 59
$ws2.Range("F26").Value = 843

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 225

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 225
$ws4.Range("F3").Value = 4
$ws4.Range("F5").Value = 5425
$ws4.Range("F6").Value = 5425
$ws4.Range("F7").Value = 191
$ws4.Range("F10").Value = 540
$ws4.Range("F11").Value = 1216
$ws4.Range("F14").Value = 779
$ws4.Range("F15").Value = 6387
$ws4.Range("F16").Value = 36
$ws4.Range("F18").Value = 115
$ws4.Range("F19").Value = 3812
$ws4.Range("F20").Value = 274
$ws4.Range("F23").Value = 4097
$ws4.Range("F24").Value = 4010
$ws4.Range("F26").Value = 194
$ws4.Range("F27").Value = 256
$ws4.Range("F32").Value = 147
$ws4.Range("F33").Value = 7170
$ws4.Range("F34").Value = 36
$ws4.Range("F35").Value = 1186
$ws4.Range("F36").Value = 582
$ws4.Range("F38").Value = 111
$ws4.Range("F39").Value = 976
$ws4.Range("F41").Value = 1446
$ws4.Range("F43").Value = 781
$ws4.Range("F45").Value = 3333
$ws4.Range("F48").Value = 801
